$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to remain text, matching original inline-string storage
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.277.64"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.749.85"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "618.64"
$ws.Range("E5").Value = "  +0.43%  "
$ws.Range("D6").Value = "181.15"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("D7").Value = "3.747.01"
$ws.Range("E7").Value = "  -1.25%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "0.534"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").Value = "6.37"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").Value = "  -3.78%  "
$ws.Range("D13").Value = "40.09"
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "0.0000255"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "4.385.45"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "3.755.45"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "70.374.19"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  -2.11%  "
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "16.52"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").Value = "504.04"
$ws.Range("E21").Value = "  -2.73%  "
$ws.Range("D22").Value = "9.20"
$ws.Range("E22").Value = "  -3.07%  "
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("D24").Value = "2.61"
$ws.Range("E24").Value = "  +5.10%  "
$ws.Range("D25").Value = "86.47"
$ws.Range("E25").Value = "  -2.80%  "
$ws.Range("D26").Value = "13.03"
$ws.Range("E26").Value = "  -4.07%  "
$ws.Range("D27").Value = "11.28"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("D28").Value = "0.0000133"
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "2.48"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").Value = "2.93"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("D32").Value = "7.90"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").Value = "30.54"
$ws.Range("E33").Value = "  -5.12%  "
$ws.Range("D34").Value = "0.114"
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").Value = "6.13"
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("D38").Value = "0.353"
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").Value = "0.139"
$ws.Range("E39").Value = "  +5.81%  "
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").Value = "  +13.48%  "
$ws.Range("D41").Value = "2.07"
$ws.Range("E41").Value = "  -4.75%  "
$ws.Range("D42").Value = "50.00"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").Value = "45.39"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "435.27"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").Value = "8.61"
$ws.Range("E45").Value = "  -2.82%  "
$ws.Range("D46").Value = "2.955.58"
$ws.Range("E46").Value = "  -5.24%  "
$ws.Range("D47").Value = "0.0363"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "27.44"
$ws.Range("E48").Value = "  -1.41%  "
$ws.Range("D50").Value = "138.27"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "2.48"
$ws.Range("E51").Value = "  +0.17%  "
